$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(0, "struggle", 2.846464770180838, -7.078682354518345, -3.06583663395473, -0.0001527163112768, -0.024892758578062, 0.0048869219608604),
    @(100, "struggle", 2.835772105625697, -7.043793456895011, -2.961106973034996, -0.0087048299610614, -0.0183259565383195, 0.0035124751739203),
    @(200, "struggle", 2.915691324642726, -7.111980744770594, -2.798882663249969, -0.0713185146450996, -0.1441642045974731, 0.0519235469400882),
    @(300, "struggle", 3.055751519543784, -7.220968450818743, -2.894419597727912, -0.09269879758358, 0.0145080499351024, -0.007177666760981),
    @(400, "struggle", 2.955562557492938, -7.134888444628034, -2.944399050303869, -0.090408056974411, -0.0394008085131645, -0.0087048299610614),
    @(500, "struggle", 2.833344757556915, -7.338198423385621, -2.509933024644852, -0.07376197725534429, -0.102472648024559, 0.0461203269660472),
    @(600, "struggle", 3.108331612178258, -7.027578847748893, -2.627250722476415, -0.1244637966156005, -0.4952589869499206, 0.25641068816185),
    @(700, "struggle", 3.771172387259348, -7.547038768018995, -3.545775217669351, 0.955545961856842, 0.418595403432846, 0.5012149214744568),
    @(800, "struggle", 2.675201347896034, -6.581202302660263, 2.435893160956232, 0.2005165219306945, 3.583182811737061, 0.4506658315658569),
    @(900, "struggle", 0.8038081441606844, -7.229396700859067, 1.553160754697669, 0.0250454749912023, 1.217149019241333, 0.2982549667358398),
    @(1000, "struggle", 1.023633156503948, -4.972137170178557, -3.916851524795791, 0.4430300295352936, 1.049771904945374, -0.152105450630188),
    @(1100, "struggle", 0.3411836709295031, -4.24002621429307, -5.904611808913096, -0.1605048477649688, 0.6899722814559937, -0.2005165219306945),
    @(1200, "struggle", -0.6227953433990485, -4.929600954055787, -6.7005569934845, 0.2364048510789871, -1.487151384353638, -0.1020144969224929),
    @(1300, "struggle", -0.2611346755709001, -5.176580480166844, -11.08990749291009, -1.693776607513428, -5.596747398376465, 1.112385630607605),
    @(1400, "struggle", 7.22266825607841, -7.398067985262175, -7.469464063644459, -1.431409955024719, -1.327410221099854, 1.870163917541504),
    @(1500, "struggle", 17.43506016050075, -11.15922947440832, 8.020842296736639, 0.08170322328805921, -2.866485118865967, 0.2145664244890213),
    @(1600, "struggle", 4.810774190085277, -6.605375971112934, 0.3875993319920126, 0.2112066596746444, -2.239890098571777, -0.6490443348884583),
    @(1700, "struggle", 3.291332300220217, -5.354068347385952, 1.089433806283133, -0.4173736870288849, -0.5337435007095337, -0.0603229440748691),
    @(1800, "struggle", 3.138020081179482, -6.129925046648298, 3.955937453678676, 0.2202169150114059, 1.032209515571594, 0.3608686327934265),
    @(1900, "struggle", 3.258425533771516, -4.712551474571227, 5.692880451679232, 2.595261096954346, 7.212944030761719, -0.4940372705459595),
    @(2000, "struggle", -3.766383392470183, -8.431432792118596, 7.765124661581846, 0.2443460971117019, 2.484389066696167, 0.5590944290161133),
    @(2100, "struggle", 1.177985225405012, -3.925881973334721, -4.60967251232692, 0.3139847218990326, 0.8843801617622375, -0.1221730485558509),
    @(2200, "struggle", 1.052592243467076, -3.166507703917357, -6.040607401302906, -1.060767531394958, 1.06214189529419, 0.0241291765123605),
    @(2300, "struggle", -1.074721106461116, -3.308761754206249, -6.224646031856537, -0.1357648074626922, -0.5329799056053162, -0.2561052441596985),
    @(2400, "struggle", -1.158144678388323, -3.9557591165815, -8.491945947919575, -0.5283984541893005, -2.495231866836548, 0.1453859210014343),
    @(2500, "struggle", 0.7903139420917948, -5.482765521321981, -6.279002168348867, -0.8752171993255615, -1.659415483474731, -0.0487165041267871),
    @(2600, "struggle", 4.010212659835815, -5.284669637680054, -2.184034883975983, -0.0806342139840126, -0.5641340613365173, 0.087353728711605),
    @(2700, "struggle", 2.293728096144541, -7.460245260170529, 0.3245020040443973, -0.1470658034086227, -0.0910189226269722, -0.2229658216238021),
    @(2800, "struggle", 1.675990547452653, -6.526311159133912, -0.6771522419793252, 0.4641048610210418, 0.3608686327934265, 0.3602577745914459),
    @(2900, "struggle", 3.116939672401965, -5.273013770580299, -3.711685695818474, 0.2412917762994766, 0.2144137024879455, -0.0186313893646001)
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
    $ws.Cells.Item($r, 8).Value = $row[7]
    $r = $r + 1
}

Write-Host "Rows written:" ($r - 2)
